# Typo errors in slides fixed
$p = $ppt.ActivePresentation

# --- Slide 7 ("Bubble Sort"): shape 69 title wrongly said "SelectionSort" ---
$s7 = $p.Slides.Item(7)
$shBubble = $s7.Shapes.Item(2)
$shBubble.Name = "Algorithm BubbleSort(A[0..n-1])…"
$tr7 = $shBubble.TextFrame.TextRange
$full7 = $tr7.Text
$idx7 = $full7.IndexOf("SelectionSort(A[0..n-1])")
$run7 = $tr7.Characters($idx7 + 1, 24)
$run7.Text = "BubbleSort(A[0..n-1])"

# --- Slide 8 ("Insertion Sort"): shape 75 grew taller + "pos-1" -> "pos - 1" ---
$s8 = $p.Slides.Item(8)
$shInsertion = $s8.Shapes.Item(2)
$shInsertion.Height = 506.3959655761719

$tr8 = $shInsertion.TextFrame.TextRange
$full8 = $tr8.Text
$idx8 = $full8.LastIndexOf("pos-1")
$run8 = $tr8.Characters($idx8, 6)
$run8.Text = " pos - 1"

# --- Slide 15 ("Numerical Problems"): shape 117 lost a stray leftover
#     click animation on a 7th (non-existent) paragraph, and grew taller ---
$s15 = $p.Slides.Item(15)
$shNumerical = $s15.Shapes.Item(2)
$shNumerical.Top = 55.20051193237305
$shNumerical.Height = 475.7345275878906

$seq15 = $s15.TimeLine.MainSequence
$seq15.Item($seq15.Count).Delete()
